$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run Info")

# Data for the three new runs appended to the "Run Info" log (rows 203-205).
$newRows = @(
    @{ Row = 203; A = 42962.441550925927; B = "rcp45"; C = 1;     D = 10000; E = 1; F = 8.8682453696637253;  G = 0;                  H = 0.3; I = 0.1; J = 4; K = 4; L = 2; M = 0.36; N = 1.5; O = 0.46; P = 4.8444000000000003 },
    @{ Row = 204; A = 42962.480162037034; B = "rcp45"; C = 1;     D = 10000; E = 2; F = 4.2437469940933035;  G = 0;                  H = 0.3; I = 0.1; J = 4; K = 4; L = 2; M = 0.36; N = 1.5; O = 0.46; P = 4.8444000000000003 },
    @{ Row = 205; A = 42962.489386574074; B = "rcp45"; C = 1;     D = 1;     E = 6; F = 39.402462387567489; G = 5.0069930069930075; H = 0.3; I = 0.1; J = 4; K = 4; L = 2; M = 0.36; N = 1.5; O = 0.46; P = 4.8444000000000003 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy the formatting (date number format etc.) from the prior row, then
    # overwrite every cell value so only the new data's format carries over.
    $ws.Cells.Item($row - 1, 1).Copy($ws.Cells.Item($row, 1))

    $ws.Cells.Item($row, 1).Value  = $r.A
    $ws.Cells.Item($row, 2).Value  = $r.B
    $ws.Cells.Item($row, 3).Value  = $r.C
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = $r.E
    $ws.Cells.Item($row, 6).Value  = $r.F
    $ws.Cells.Item($row, 7).Value  = $r.G
    $ws.Cells.Item($row, 8).Value  = $r.H
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
}

# Match the selection state left behind by the edit (last row, all columns).
$ws.Range("A205:P205").Select() | Out-Null
